# Reorders the goods-name labels in column A for the rows whose shared
# string entries were reshuffled in the authoring tool (see commit
# "added csv output for goods aggregations"). The numeric counts in
# column B stay attached to their original row/position; only the
# label text shown in column A moves, mirroring the shared-strings
# table reshuffle from the source diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A15").Value = "мелочь"
$ws.Range("A16").Value = "особливый товар"

$ws.Range("A17").Value = "серебреный товар"
$ws.Range("A18").Value = "деревенский товар"

$ws.Range("A19").Value = "крамными товар"
$ws.Range("A20").Value = "небогатый товар"

$ws.Range("A24").Value = "щепетильный товар"
$ws.Range("A25").Value = "нужный товар"

$ws.Range("A29").Value = "медный товар"
$ws.Range("A31").Value = "питейный припасы"

$ws.Range("A35").Value = "произрастание"
$ws.Range("A36").Value = "заморский товар"
$ws.Range("A37").Value = "галантерейный товар"
$ws.Range("A38").Value = "купецкий товар"

$ws.Range("A39").Value = "надлежащий товар"
$ws.Range("A40").Value = "домовый товар"

$ws.Range("A41").Value = "харчевой припасы"
$ws.Range("A42").Value = "рукодельный товар"
